$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would be mis-parsed as a number by Excel (e.g. "1.000" -> 1,
# or "0.00000000120" -> 1.2E-09). Force these to remain literal text via a text number
# format before assigning the value (classic Excel "store as text" behavior).
$textForceCells = @("D4", "D5", "D6", "D8", "D11", "D13", "D14", "D15", "D16", "D18", "D20", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D33", "D34", "D36", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply new values
$ws.Range('D4').Value = '0.9995'
$ws.Range('D5').Value = '240.04'
$ws.Range('D6').Value = '0.6298'
$ws.Range('D8').Value = '0.07626'
$ws.Range('D11').Value = '0.07737'
$ws.Range('D13').Value = '0.00001119'
$ws.Range('D14').Value = '5.002'
$ws.Range('D15').Value = '0.6783'
$ws.Range('D16').Value = '83.69'
$ws.Range('D18').Value = '6.183'
$ws.Range('D20').Value = '228.74'
$ws.Range('D23').Value = '7.485'
$ws.Range('D24').Value = '1.000'
$ws.Range('D25').Value = '157.15'
$ws.Range('D26').Value = '0.1396'
$ws.Range('D27').Value = '8.345'
$ws.Range('D30').Value = '1.301'
$ws.Range('D31').Value = '0.05590'
$ws.Range('D33').Value = '4.029'
$ws.Range('D34').Value = '1.848'
$ws.Range('D36').Value = '0.7107'
$ws.Range('D39').Value = '0.01805'
$ws.Range('D40').Value = '2.774'
$ws.Range('D41').Value = '6.398'
$ws.Range('D42').Value = '0.9049'
$ws.Range('D43').Value = '1.000'
$ws.Range('D44').Value = '101.98'
$ws.Range('D45').Value = '65.91'
$ws.Range('D46').Value = '0.00000000120'
$ws.Range('D47').Value = '7.130'
$ws.Range('D48').Value = '0.4014'
$ws.Range('D49').Value = '9.031'
$ws.Range('D50').Value = '1.683'

# Plain value updates (already unambiguous text)
$ws.Range('D2').Value = '29.381.61'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '1.847.73'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E5').Value = '  -0.37%  '
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('E10').Value = '  -1.02%  '
$ws.Range('E11').Value = '  +0.03%  '
$ws.Range('D12').Value = '1.852.80'
$ws.Range('E12').Value = '  -6.66%  '
$ws.Range('E13').Value = '  +12.36%  '
$ws.Range('E14').Value = '  +0.27%  '
$ws.Range('E15').Value = '  -0.49%  '
$ws.Range('E16').Value = '  +0.85%  '
$ws.Range('D17').Value = '2.106.02'
$ws.Range('E17').Value = '  -6.99%  '
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = '29.401.74'
$ws.Range('E19').Value = '  +0.00%  '
$ws.Range('E20').Value = '  -0.52%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('E23').Value = '  -0.76%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('E26').Value = '  -0.15%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  +3.88%  '
$ws.Range('E31').Value = '  -2.04%  '
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('E33').Value = '  +0.50%  '
$ws.Range('E34').Value = '  +0.38%  '
$ws.Range('E35').Value = '  +0.20%  '
$ws.Range('E36').Value = '  -0.66%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('D38').Value = '1.239.40'
$ws.Range('E38').Value = '  -0.15%  '
$ws.Range('E39').Value = '  -0.17%  '
$ws.Range('E40').Value = '  -0.89%  '
$ws.Range('E41').Value = '  +5.01%  '
$ws.Range('E42').Value = '  +0.51%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E47').Value = '  +1.23%  '
$ws.Range('B48').Value = 'TheSandbox'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E48').Value = '  +0.06%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E49').Value = '  -0.95%  '
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('E51').Value = '  -0.04%  '
